# Auto commit: 2025-04-22 12:51:12
# Rename the four "Day" statistics sheets (drop the "Day" token) and
# refresh the APPLE (column F) statistics with their recalculated values.

$wb = $excel.ActiveWorkbook

# --- 1) Rename worksheets ------------------------------------------------
$wb.Worksheets.Item("1DDay_Statistics").Name  = "1D_Statistics"
$wb.Worksheets.Item("3DDay_Statistics").Name  = "3D_Statistics"
$wb.Worksheets.Item("5DDay_Statistics").Name  = "5D_Statistics"
$wb.Worksheets.Item("10DDay_Statistics").Name = "10D_Statistics"

# --- 2) Update column F (APPLE) statistics on every sheet ---------------
$ws1 = $wb.Worksheets.Item("1D_Statistics")
$ws1.Range("F3").Value  = 0.001754828878540432
$ws1.Range("F4").Value  = 0.01254982721302963
$ws1.Range("F5").Value  = 0.000157498163076899
$ws1.Range("F6").Value  = -0.05249309071118011
$ws1.Range("F7").Value  = 1.470640435624273
$ws1.Range("F8").Value  = -0.01732573067276034
$ws1.Range("F9").Value  = -0.02534237462518862
$ws1.Range("F10").Value = -0.02841681823381926
$ws1.Range("F11").Value = -0.03844888055246375

$ws2 = $wb.Worksheets.Item("3D_Statistics")
$ws2.Range("F3").Value  = 0.005264486635621296
$ws2.Range("F4").Value  = 0.02173693835917784
$ws2.Range("F5").Value  = 0.0004724944892306971
$ws2.Range("F6").Value  = -0.05249309071118011
$ws2.Range("F7").Value  = 1.470640435624273
$ws2.Range("F8").Value  = -0.03000904580347542
$ws2.Range("F9").Value  = -0.04389428043527098
$ws2.Range("F10").Value = -0.04921937297042464
$ws2.Range("F11").Value = -0.06659541461101413

$ws3 = $wb.Worksheets.Item("5D_Statistics")
$ws3.Range("F3").Value  = 0.00877414439270216
$ws3.Range("F4").Value  = 0.02806226675421099
$ws3.Range("F5").Value  = 0.0007874908153844952
$ws3.Range("F6").Value  = -0.05249309071118011
$ws3.Range("F7").Value  = 1.470640435624273
$ws3.Range("F8").Value  = -0.03874151154414529
$ws3.Range("F9").Value  = -0.05666727237318752
$ws3.Range("F10").Value = -0.06354193727507539
$ws3.Range("F11").Value = -0.08597431057407863

$ws4 = $wb.Worksheets.Item("10D_Statistics")
$ws4.Range("F3").Value  = 0.01754828878540432
$ws4.Range("F4").Value  = 0.03968603823473679
$ws4.Range("F5").Value  = 0.00157498163076899
$ws4.Range("F6").Value  = -0.05249309071118011
$ws4.Range("F7").Value  = 1.470640435624273
$ws4.Range("F8").Value  = -0.0547887710525641
$ws4.Range("F9").Value  = -0.080139625132852
$ws4.Range("F10").Value = -0.08986186947387211
$ws4.Range("F11").Value = -0.1215860360295386

$ws5 = $wb.Worksheets.Item("Annual_Statistics")
$ws5.Range("F3").Value  = 0.438707219635108
$ws5.Range("F4").Value  = 0.198430191173684
$ws5.Range("F5").Value  = 0.03937454076922476
$ws5.Range("F6").Value  = -0.05249309071118011
$ws5.Range("F7").Value  = 1.470640435624273
$ws5.Range("F8").Value  = -0.2739438552628204
$ws5.Range("F9").Value  = -0.4006981256642599
$ws5.Range("F10").Value = -0.4493093473693605
$ws5.Range("F11").Value = -0.6079301801476928
